# RS-RIG schema doc: update the "type" field row, insert a new
# "vehiculeType" row after it, fill in the "plate" row's description,
# and fill in the "coord" row's format cell.

function Get-TargetTable($doc) {
    # The same "Nom de balise / Champ correspondant / ..." header is reused
    # by several tables in this document, so locate the one table that
    # actually contains the "type" / "Type" field row we need to edit
    # (cell 1 == "type" and cell 2 == "Type"), rather than trusting a
    # hard-coded table index.
    for ($ti = 1; $ti -le $doc.Tables.Count; $ti++) {
        $tbl = $doc.Tables.Item($ti)
        if ($tbl.Columns.Count -ne 6) { continue }
        for ($ri = 1; $ri -le $tbl.Rows.Count; $ri++) {
            $c1 = $tbl.Cell($ri, 1).Range.Text.TrimEnd([char]7, [char]13)
            $c2 = $tbl.Cell($ri, 2).Range.Text.TrimEnd([char]7, [char]13)
            if ($c1 -eq "type" -and $c2 -eq "Type") {
                return $tbl
            }
        }
    }
    return $null
}

function Get-RowIndexByFirstCell($table, $text) {
    for ($i = 1; $i -le $table.Rows.Count; $i++) {
        $c = $table.Cell($i, 1).Range.Text.TrimEnd([char]7, [char]13)
        if ($c -eq $text) {
            return $i
        }
    }
    return -1
}

$d = $word.ActiveDocument
$t = Get-TargetTable $d

# --- 1) "type" row: rename tag + field label, shorten description ---
$typeRow = Get-RowIndexByFirstCell $t "type"

$t.Cell($typeRow, 1).Range.Text = "resourceType"
$t.Cell($typeRow, 2).Range.Text = "Type de ressource"
$t.Cell($typeRow, 5).Range.Text = "Type de ressource mobilisée (type moyen)"

# --- 2) insert a new row right after the (now renamed) "type" row ---
$plateRowBefore = Get-RowIndexByFirstCell $t "plate"
$newRow = $t.Rows.Add($t.Rows.Item($plateRowBefore))
$newRowIndex = $newRow.Index

$t.Cell($newRowIndex, 1).Range.Text = "vehiculeType"
$t.Cell($newRowIndex, 2).Range.Text = "Type de vecteur"
$t.Cell($newRowIndex, 3).Range.Text = "string"
$t.Cell($newRowIndex, 4).Range.Text = "0..1"
$t.Cell($newRowIndex, 5).Range.Text = "Type de ressource mobilisée (nomenclature type de vecteur à implémenter)"

# --- 3) "plate" row: fill in the previously-empty description cell ---
$plateRow = Get-RowIndexByFirstCell $t "plate"
$t.Cell($plateRow, 5).Range.Text = "N° d'immatriculation du vecteur"

# --- 4) "coord" row: fill in the previously-empty format cell ---
$coordRow = Get-RowIndexByFirstCell $t "coord"
$t.Cell($coordRow, 3).Range.Text = "coord"
